# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy style (bold, centered, bordered) from the last existing
# header cell (AC1) into the three new header cells, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-42: every player row gets the same team record for 2005.
$lastRow = 42
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
